$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 32) to the "Horas trabajadas" log.
$ws.Range("A32").Value = "Bruno Díaz"
$ws.Range("B31").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B32").Value = 42862
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = "Sprint 1 - MER"
$ws.Range("E32").Value = "Modificando el MER a partir de la lista de requerimientos"

# Update the active selection to match the authored state.
[void]$ws.Range("E26").Select()
